$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same rows of data that were
# updated in the source (ticket/attendance counts refreshed on regeneration).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1572
    $ws.Range("F5").Value = 8
    $ws.Range("F10").Value = 429
}
